$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95, shifting existing rows 95:129 down to 96:130.
$ws.Rows("95:95").Insert()

# Populate the newly inserted row 95 with the new price-report record.
$ws.Cells.Item(95, 1).Value = 7
$ws.Cells.Item(95, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(95, 3).Value = "Ñuble"
$ws.Cells.Item(95, 4).Value = 45027
$ws.Cells.Item(95, 5).Value = 16
$ws.Cells.Item(95, 6).Value = 100112031
$ws.Cells.Item(95, 7).Value = "Poroto verde"
$ws.Cells.Item(95, 8).Value = "Sin especificar"
$ws.Cells.Item(95, 9).Value = "Primera"
$ws.Cells.Item(95, 10).Value = 80
$ws.Cells.Item(95, 11).Value = 27000
$ws.Cells.Item(95, 12).Value = 28000
$ws.Cells.Item(95, 13).Value = 27500
$ws.Cells.Item(95, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(95, 15).Value = "Región del Maule"
$ws.Cells.Item(95, 16).Value = 1100
$ws.Cells.Item(95, 17).Value = 25
$ws.Cells.Item(95, 18).Value = "Hortaliza"
